$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" - refresh the localization-status report with
# a new handoff run: new source/target file ids, new timestamps, a "Ready
# for handoff" status (instead of the old handed-back status), and removal
# of the now-stale "Latest Target File" / "Latest Handback File" columns
# (F/G) on the per-locale sheets, since nothing has been handed back yet.
# ---------------------------------------------------------------------------

$oldFile1 = "31bf859f-1ce4-4a73-8d74-02a1f0514d03.md"
$oldFile2 = "f0965c0a-12b0-4588-a577-42a536641f9f.md"
$newFile1 = "1dac9d35-f631-48f4-9714-353f7d407fc7.md"
$newFile2 = "ffff30bd7555-9088-4152-949e-1fdcf6947d9b.md"

$oldXlf1zh = "31bf859f-1ce4-4a73-8d74-02a1f0514d03.4b0e2bb43cdf2bf17d4649f7ef01886115e27aab.zh-cn.xlf"
$oldXlf2zh = "f0965c0a-12b0-4588-a577-42a536641f9f.ac496932139a0ef0078454db283a05d5d94f1486.zh-cn.xlf"
$newXlfzh  = "1dac9d35-f631-48f4-9714-353f7d407fc7.50d3ccbbe921dc5ccb5afe7d33f9650c9dd2e724.zh-cn.xlf"

$oldXlf1de = "31bf859f-1ce4-4a73-8d74-02a1f0514d03.4b0e2bb43cdf2bf17d4649f7ef01886115e27aab.de-de.xlf"
$oldXlf2de = "f0965c0a-12b0-4588-a577-42a536641f9f.ac496932139a0ef0078454db283a05d5d94f1486.de-de.xlf"
$newXlfde  = "1dac9d35-f631-48f4-9714-353f7d407fc7.50d3ccbbe921dc5ccb5afe7d33f9650c9dd2e724.de-de.xlf"

$newStatus = "Ready for handoff"
$newHandoffDate = "2016-03-24 21:23:42"
$newHandoffDatetime = "2016-03-24 21:23:37"
$newHandbackDatetime = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("D2").Value = $newHandoffDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $newHandoffDate

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("D2").Value = $newXlfzh
$wsZh.Range("E2").Value = $newHandoffDatetime
$wsZh.Range("H2").Value = $newHandbackDatetime
$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Clear()

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("D3").Value = $newXlfzh
$wsZh.Range("E3").Value = $newHandoffDatetime
$wsZh.Range("H3").Value = $newHandbackDatetime
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Clear()

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("D2").Value = $newXlfde
$wsDe.Range("H2").Value = $newHandbackDatetime
$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Clear()

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("D3").Value = $newXlfde
$wsDe.Range("H3").Value = $newHandbackDatetime
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Clear()
